$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Brainstorm results: fill in the blank "Requisito" names that were
#     already sketched out in rows 13-25, then append many more new
#     requirement rows gathered during the brainstorm session. ---

# Row 26 used to be the tables thick-bottom-border closing row; drop it
# entirely so the left-over heavy border / row height do not linger once
# it becomes an ordinary data row again.
$ws.Rows.Item(26).Delete()

# A13:A19, A21:A25 already carry the correct body-row cell style - just
# fill in the text. A20 is missing a cell entirely in the source sheet,
# so clone the style from its neighbour before writing into it.
$ws.Range("A13").Value = "Terminal mobile de atendimento;"
$ws.Range("A14").Value = "Sensor de vagas livres/ocupadas"
$ws.Range("A15").Value = "Controle de filas de entrada e saída"
$ws.Range("A16").Value = "Controle de tickets"
$ws.Range("A17").Value = "Tratamento de meios de pagamento"
$ws.Range("A18").Value = "Tipificação de uso (mensalista/avulso/conveniado)"
$ws.Range("A19").Value = "Administração de caixa"
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A20").Value = "Configuração de tabelas preços (dias/horários/tipo de uso)"
$ws.Range("A21").Value = "Emissão de nota fiscal"
$ws.Range("A22").Value = "Gestão de cobranças de mensalistas e conveniados"
$ws.Range("A23").Value = "Reconhecimento de placa para recuperar dados cadastrais"
$ws.Range("A24").Value = "Cadastro de veículos e clientes"
$ws.Range("A25").Value = "Reservar vagas"

# Build 15 fresh body rows (26-40) by cloning the formatting already used
# by the other requirement rows (fills/borders on A:E), then fill column A.
$ws.Range("A19:E19").Copy()
$ws.Range("A26:E40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A26").Value = "Solução mobile para o cliente fazer reservas e pagamentos"
$ws.Range("A27").Value = "Orientação por voz"
$ws.Range("A28").Value = "Integração com TAG"
$ws.Range("A29").Value = "Relatório de gerenciamento de média de ocupação de vagas por dia e horário"
$ws.Range("A30").Value = "Relatório de gerenciamento do valor recebido por meio de pagamento por mês"
$ws.Range("A31").Value = "Forum/Sistema de Denuncias"
$ws.Range("A32").Value = "Estacionamento Vertical com elevador"
$ws.Range("A33").Value = "Planta digital"
$ws.Range("A34").Value = "Sistema de fidelidade"
$ws.Range("A35").Value = "Controle do período estacionado"
$ws.Range("A36").Value = "Sistema de vigilância automatizado, com câmeras inteligentes. "
$ws.Range("A37").Value = "Ticket digital via e-mail ou sms"
$ws.Range("A38").Value = "Sistema de Segurança contra roubos "
$ws.Range("A40").Value = "Sistema de redirecionamento para o condutor indicando vagas livres"
$ws.Range("A39").Value = "Mostrar quantidade de vagas disponíveis em determinado local"

# Column A/E got a bit narrower/wider once the longer brainstorm text was
# in place; match the resized widths from the saved file.
$ws.Columns.Item(1).ColumnWidth = 73.877604
$ws.Columns.Item(5).ColumnWidth = 86.022135

# The author scrolled/zoomed down to the newly added rows before saving.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 28
$win.ScrollColumn = 1
[void]$ws.Range("E53").Select()
